# Adds two new columns, I (header "I0") and J (header "IF"), to Sheet1,
# filling in values for rows 2 through 35, matching the existing H-column
# ("IP") header style for the new header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header cell (H1) onto the two
# new header cells so they pick up the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Set the new header labels.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: row number, I0 value, IF value.
$data = @(
    @(2, 1, 5),
    @(3, 6, 8),
    @(4, 7, 8),
    @(5, 8, 9),
    @(6, 6, 7),
    @(7, 3, 7),
    @(8, 2, 6),
    @(9, 2, 7),
    @(10, 4, 8),
    @(11, 4, 7),
    @(12, 3, 7),
    @(13, 5, 6),
    @(14, 4, 7),
    @(15, 1, 3),
    @(16, 1, 4),
    @(17, 1, 6),
    @(18, 1, 5),
    @(19, 1, 5),
    @(20, 1, 4),
    @(21, 1, 1),
    @(22, 1, 5),
    @(23, 1, 5),
    @(24, 1, 2),
    @(25, 1, 6),
    @(26, 1, 5),
    @(27, 1, 3),
    @(28, 1, 6),
    @(29, 1, 6),
    @(30, 1, 4),
    @(31, 1, 5),
    @(32, 1, 5),
    @(33, 1, 4),
    @(34, 1, 3),
    @(35, 1, 2)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $i0 = $entry[1]
    $if = $entry[2]
    $ws.Cells.Item($r, 9).Value = $i0
    $ws.Cells.Item($r, 10).Value = $if
}

Write-Host "I0 and IF columns added"
